$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 50, shifting the existing rows
# 50-128 down to 51-129 (this also extends the sheet dimension to
# A1:R129 and keeps formatting, e.g. the date style on column D, intact).
$ws.Rows(50).Insert()

# Populate the newly inserted row 50 with the new weekly record.
$ws.Range("A50").Value = 1
$ws.Range("B50").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C50").Value = "Arica y Parinacota"
$ws.Range("D50").Value = 44930
$ws.Range("E50").Value = 15
$ws.Range("F50").Value = 100114001
$ws.Range("G50").Value = "Papa"
$ws.Range("H50").Value = "Cardinal"
$ws.Range("I50").Value = "1a (cosecha)"
$ws.Range("J50").Value = 900
$ws.Range("K50").Value = 14000
$ws.Range("L50").Value = 15000
$ws.Range("M50").Value = 14333
$ws.Range("N50").Value = "$/saco 25 kilos"
$ws.Range("O50").Value = "Región de Coquimbo"
$ws.Range("P50").Value = 573
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
